# Update cryptocurrency price/volume data to latest scrape.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.358.55'
$ws.Cells.Item(2, 5).Value = '  +0.04%  '
$ws.Cells.Item(3, 4).Value = '3.321.85'
$ws.Cells.Item(3, 5).Value = '  +1.09%  '
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).Value = '''581.75'
$ws.Cells.Item(5, 5).Value = '  -0.02%  '
$ws.Cells.Item(6, 4).Value = '''175.97'
$ws.Cells.Item(6, 5).Value = '  -3.08%  '
$ws.Cells.Item(7, 4).Value = '''1.00'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '
$ws.Cells.Item(8, 4).Value = '''0.589'
$ws.Cells.Item(8, 5).Value = '  +0.49%  '
$ws.Cells.Item(9, 4).Value = '3.316.87'
$ws.Cells.Item(9, 5).Value = '  +1.08%  '
$ws.Cells.Item(10, 5).Value = '  +0.28%  '
$ws.Cells.Item(11, 5).Value = '  +0.22%  '
$ws.Cells.Item(12, 4).Value = '''45.45'
$ws.Cells.Item(12, 5).Value = '  -1.50%  '
$ws.Cells.Item(13, 5).Value = '  -1.59%  '
$ws.Cells.Item(14, 4).Value = '''658.70'
$ws.Cells.Item(14, 5).Value = '  +3.61%  '
$ws.Cells.Item(15, 4).Value = '3.862.99'
$ws.Cells.Item(15, 5).Value = '  +1.15%  '
$ws.Cells.Item(16, 5).Value = '  +0.32%  '
$ws.Cells.Item(17, 4).Value = '67.537.67'
$ws.Cells.Item(17, 5).Value = '  +0.00%  '
$ws.Cells.Item(18, 5).Value = '  -0.23%  '
$ws.Cells.Item(19, 4).Value = '3.323.57'
$ws.Cells.Item(19, 5).Value = '  +0.62%  '
$ws.Cells.Item(20, 4).Value = '''17.37'
$ws.Cells.Item(20, 5).Value = '  -0.95%  '
$ws.Cells.Item(21, 4).Value = '''10.94'
$ws.Cells.Item(21, 5).Value = '  +0.74%  '
$ws.Cells.Item(22, 4).Value = '''0.888'
$ws.Cells.Item(22, 5).Value = '  -0.39%  '
$ws.Cells.Item(23, 5).Value = '  +8.84%  '
$ws.Cells.Item(24, 4).Value = '''17.04'
$ws.Cells.Item(24, 5).Value = '  -3.33%  '
$ws.Cells.Item(25, 4).Value = '''99.50'
$ws.Cells.Item(25, 5).Value = '  +2.24%  '
$ws.Cells.Item(26, 4).Value = '''3.85'
$ws.Cells.Item(26, 5).Value = '  -3.15%  '
$ws.Cells.Item(27, 5).Value = '  -4.29%  '
$ws.Cells.Item(28, 4).Value = '''9.26'
$ws.Cells.Item(28, 5).Value = '  -2.71%  '
$ws.Cells.Item(29, 4).Value = '''33.66'
$ws.Cells.Item(29, 5).Value = '  +3.34%  '
$ws.Cells.Item(30, 4).Value = '''7.47'
$ws.Cells.Item(30, 5).Value = '  +12.48%  '
$ws.Cells.Item(31, 5).Value = '  -0.86%  '
$ws.Cells.Item(32, 4).Value = '''573.22'
$ws.Cells.Item(32, 5).Value = '  -2.81%  '
$ws.Cells.Item(33, 5).Value = '  +0.83%  '
$ws.Cells.Item(34, 5).Value = '  +0.57%  '
$ws.Cells.Item(35, 5).Value = '  +0.18%  '
$ws.Cells.Item(36, 4).Value = '3.691.29'
$ws.Cells.Item(36, 5).Value = '  -5.99%  '
$ws.Cells.Item(37, 4).Value = '''56.49'
$ws.Cells.Item(37, 5).Value = '  +1.53%  '
$ws.Cells.Item(38, 4).Value = '''3.36'
$ws.Cells.Item(38, 5).Value = '  -6.10%  '
$ws.Cells.Item(39, 4).Value = '''34.48'
$ws.Cells.Item(39, 5).Value = '  +5.28%  '
$ws.Cells.Item(40, 4).Value = '''0.130'
$ws.Cells.Item(40, 5).Value = '  +2.05%  '
$ws.Cells.Item(41, 5).Value = '  -2.39%  '
$ws.Cells.Item(42, 5).Value = '  -3.98%  '
$ws.Cells.Item(43, 5).Value = '  -0.32%  '
$ws.Cells.Item(44, 4).Value = '0.0₃0667'
$ws.Cells.Item(44, 5).Value = '  -2.32%  '
$ws.Cells.Item(45, 5).Value = '  -2.11%  '
$ws.Cells.Item(46, 5).Value = '  -1.67%  '
$ws.Cells.Item(47, 2).Value = 'Stellar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(47, 4).Value = '''0.128'
$ws.Cells.Item(47, 5).Value = '  +0.09%  '
$ws.Cells.Item(48, 2).Value = 'ThetaToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(48, 4).Value = '''2.58'
$ws.Cells.Item(48, 5).Value = '  +2.17%  '
$ws.Cells.Item(49, 5).Value = '  -0.35%  '
$ws.Cells.Item(50, 5).Value = '  +1.66%  '
$ws.Cells.Item(51, 4).Value = '''127.45'
$ws.Cells.Item(51, 5).Value = '  -2.26%  '
